$p = $ppt.ActivePresentation

# 1) Slide 6: table's style id {1B8E2D73-...} -> {819441FD-...}
$s = $p.Slides.Item(6)
$tblShape = $s.Shapes.Item(2)
$tblShape.Table.ApplyStyle("{819441FD-3CCE-40E2-B018-59A5F0F7509F}")

# 2) Swap the presentation's theme color scheme from "Integral" to the
#    stock "Office" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$tcs = $p.Slides.Item(1).ThemeColorScheme
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
